$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date for each data row (rows 2-150).
# Bump every one of those dates forward by a single day (2023-09-11 -> 2023-09-12).
$range = $ws.Range("C2:C150")
foreach ($cell in $range.Cells) {
    $currentValue = $cell.Value()
    $cell.Value = $currentValue.AddDays(1)
}
